# Add parsing raw tables (test_importtables.py etc.)
# Appends two new "tbl2Procedure" error rows (520, 530) to the Errors_ sheet,
# matching the style used for the existing tbl1Procedure error block but with
# its own fill color.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New fill color for the tbl2Procedure section (Gold, Accent 4, Lighter 60%)
$rng = $ws.Range("A37:D38")
$rng.Interior.Color = 10086143

# Row 37: header/"Base" row for the tbl2Procedure checks
$ws.Range("A37").Value = 520
$ws.Range("B37").Value = "CheckDataFrame"
$ws.Range("C37").Value = "tbl2Procedure"
$ws.Range("D37").Value = "Base"

# Row 38: new tbl2 regex-format error message
$ws.Range("A38").Value = 530
$ws.Range("B38").Value = "CheckDataFrame"
$ws.Range("C38").Value = "tbl2Procedure"
$ws.Range("D38").Value = "ERROR: tbl2 specified column contains values that don't meet required  regex format"

[void]$ws.Range("A38").Select()
